$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural shape: mirrors the row insert/delete pattern of the real edit ---
# Insert a new row for "Holberg Suite" right after the header
$ws.Rows.Item(2).Insert()

# Remove two placeholder rows (their slots get re-used further down)
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(9).Delete()

# Insert two new rows at the bottom for "Waltz No. 2" and "Winter Melodies"
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()

# --- Fill in row 2: Holberg Suite / Edward Grieg ---
$ws.Cells.Item(2,1).Value = "Holberg Suite"
$ws.Cells.Item(2,2).Value = "Edward Grieg"
$ws.Cells.Item(2,3).Value = ""
$ws.Cells.Item(2,4).Value = ""
$ws.Cells.Item(4,5).Copy()
$ws.Cells.Item(2,5).PasteSpecial()
$ws.Cells.Item(2,6).Value = "Mvmts 1, 3, 5 ONLY"

# --- Update row 8: The Elements / Richard Meyer ---
$ws.Cells.Item(8,1).Value = "The Elements"
$ws.Cells.Item(8,2).Value = "Richard Meyer"
$ws.Cells.Item(8,3).Value = ""
$ws.Cells.Item(8,4).Value = "Hal Leonard"
$ws.Cells.Item(4,5).Copy()
$ws.Cells.Item(8,5).PasteSpecial()
$ws.Cells.Item(8,6).Value = "Really cool piece!"

# --- Fill in row 11: Waltz No. 2 / Shostakovich ---
$ws.Cells.Item(11,1).Value = "Waltz No. 2"
$ws.Cells.Item(11,2).Value = "Shostakovich"
$ws.Cells.Item(11,3).Value = "TEst"
$ws.Cells.Item(11,4).Value = "Test"
$ws.Cells.Item(11,5).Formula = '="4"'
$ws.Cells.Item(11,5).Copy()
$ws.Cells.Item(11,5).PasteSpecial(-4163)
$ws.Cells.Item(11,6).Value = "Cello solo!"

# --- Fill in row 12: Winter Melodies / Reese, M. ---
$ws.Cells.Item(12,1).Value = "Winter Melodies"
$ws.Cells.Item(12,2).Value = "Reese, M."
$ws.Cells.Item(12,3).Value = ""
$ws.Cells.Item(12,4).Value = "Reese Music"
$ws.Cells.Item(4,5).Copy()
$ws.Cells.Item(12,5).PasteSpecial()
$ws.Cells.Item(12,6).Value = "This is a test."
